# Applies the "automatic update of files" change:
#  - Column C ("Förändrad") is refreshed to the new timestamp serial 46072
#    for every data row (2..33).
#  - Rows 10..33 represent a rolling log of cases; each row's identifying
#    data (Beteckning, Datum, Markägare, Area) shifts down into the next
#    row, with the last row's data wrapping around into the first row of
#    that block (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 33
$newChanged   = 46072

# 1) Refresh column C ("Förändrad") for every data row.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newChanged
}

# 2) Rotate the Beteckning (A), Datum (B), Markägare (F) and Area (G)
#    columns for rows 10..33 down by one row, wrapping the last row's
#    values around to the first row of the block.
$blockFirst = 10
$blockLast  = 33

$colA = @{}
$colB = @{}
$colF = @{}
$colG = @{}

for ($r = $blockFirst; $r -le $blockLast; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
    $colF[$r] = $ws.Cells.Item($r, 6).Value2
    $colG[$r] = $ws.Cells.Item($r, 7).Value2
}

for ($r = $blockFirst; $r -le $blockLast; $r++) {
    if ($r -eq $blockFirst) {
        $srcRow = $blockLast
    } else {
        $srcRow = $r - 1
    }

    $ws.Cells.Item($r, 1).Value = $colA[$srcRow]
    $ws.Cells.Item($r, 2).Value = $colB[$srcRow]
    $ws.Cells.Item($r, 6).Value = $colF[$srcRow]
    $ws.Cells.Item($r, 7).Value = $colG[$srcRow]
}
